{"js": "// Replace the date line and every two-digit-by-two-digit multiplication\n// answer in the table with the new values from the commit.\n// Every non-empty paragraph in the document body (the title line plus the\n// 25 filled-in table cells) gets exactly one new value, in document order.\nconst replacements = [\n  \"2025-07-17 Thursday\",\n  \"72\u00d756=4032\",\n  \"54\u00d731=1674\",\n  \"39\u00d734=1326\",\n  \"41\u00d788=3608\",\n  \"60\u00d753=3180\",\n  \"98\u00d713=1274\",\n  \"82\u00d729=2378\",\n  \"91\u00d791=8281\",\n  \"80\u00d715=1200\",\n  \"65\u00d769=4485\",\n  \"35\u00d797=3395\",\n  \"15\u00d741=615\",\n  \"54\u00d730=1620\",\n  \"66\u00d733=2178\",\n  \"73\u00d733=2409\",\n  \"94\u00d769=6486\",\n  \"39\u00d770=2730\",\n  \"16\u00d749=784\",\n  \"49\u00d769=3381\",\n  \"86\u00d734=2924\",\n  \"77\u00d755=4235\",\n  \"93\u00d737=3441\",\n  \"53\u00d768=3604\",\n  \"13\u00d746=598\",\n  \"31\u00d790=2790\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst nonEmpty = paragraphs.items.filter((p) => p.text && p.text.length > 0);\n\nif (nonEmpty.length !== replacements.length) {\n  throw new Error(\n    `Expected ${replacements.length} non-empty paragraphs, found ${nonEmpty.length}`\n  );\n}\n\nfor (let i = 0; i < nonEmpty.length; i++) {\n  nonEmpty[i].insertText(replacements[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and every two-digit-by-two-digit multiplication\n# answer in the table with the new values from the commit.\n# Every non-empty paragraph in the document (the title line plus the\n# 25 filled-in table cells) gets exactly one new value, in document order.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    \"2025-07-17 Thursday\",\n    \"72\u00d756=4032\",\n    \"54\u00d731=1674\",\n    \"39\u00d734=1326\",\n    \"41\u00d788=3608\",\n    \"60\u00d753=3180\",\n    \"98\u00d713=1274\",\n    \"82\u00d729=2378\",\n    \"91\u00d791=8281\",\n    \"80\u00d715=1200\",\n    \"65\u00d769=4485\",\n    \"35\u00d797=3395\",\n    \"15\u00d741=615\",\n    \"54\u00d730=1620\",\n    \"66\u00d733=2178\",\n    \"73\u00d733=2409\",\n    \"94\u00d769=6486\",\n    \"39\u00d770=2730\",\n    \"16\u00d749=784\",\n    \"49\u00d769=3381\",\n    \"86\u00d734=2924\",\n    \"77\u00d755=4235\",\n    \"93\u00d737=3441\",\n    \"53\u00d768=3604\",\n    \"13\u00d746=598\",\n    \"31\u00d790=2790\"\n)\n\n# Snapshot the paragraph objects first so mutating their text while we\n# iterate doesn't perturb the live collection/indexing.\n$paras = @()\nforeach ($p in $d.Paragraphs) {\n    $paras += $p\n}\n\n$idx = 0\nforeach ($p in $paras) {\n    $text = $p.Range.Text\n    $trimmed = $text.Trim([char]13, [char]7)\n    if ($trimmed.Length -gt 0) {\n        $p.Range.Text = $replacements[$idx]\n        $idx++\n    }\n}\n\nif ($idx -ne $replacements.Count) {\n    throw \"Expected $($replacements.Count) non-empty paragraphs, replaced $idx\"\n}\n"}
